$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.15869999999999
$ws.Range("C3").Value = -11.812
$ws.Range("C5").Value = -12.5625
$ws.Range("B9").Value = 8.668900000000002
$ws.Range("C11").Value = -13.2797
$ws.Range("C12").Value = -14.47700000000002
$ws.Range("B13").Value = 5.496700000000003
$ws.Range("B16").Value = 9.237500000000006
$ws.Range("B18").Value = 5.268700000000003
$ws.Range("B20").Value = 5.777799999999997
$ws.Range("C21").Value = -12.80580000000001
